$data = @(
    @(13,1.080166504188284,-1.390476677927515,0.2980841856143616,0.471029234077814,0.3747115521091721),
    @(14,0.8994690427722825,-1.398574963125803,0.2952507445860527,0.4712791868485272,0.3736449071541705),
    @(15,0.4261995378958547,-1.421233767473362,0.2901329598009249,0.4707988478433893,0.3687079851421951),
    @(16,-0.1293866511906445,-1.445894254428907,0.2881183117335868,0.4634287222226601,0.3671734660080277),
    @(17,0.2056862336605694,-1.431506400279586,0.2916483634698089,0.4621820912026547,0.3668967082701994),
    @(18,-0.1867035728334181,-1.449493393027206,0.2920552357523991,0.4634912525595197,0.3664404882586997),
    @(19,-0.6440158654643769,-1.469374742207627,0.2927357270699735,0.476967520816862,0.3660762421732142),
    @(20,-1.00310434740966,-1.486577319749939,0.2823257794302023,0.4573051486249445,0.3609390017017186),
    @(21,-1.504651547717967,-1.504352060833916,0.2824063043884478,0.4673525213334017,0.3589114020489111),
    @(22,-1.523057165276164,-1.506251585126765,0.2874913509468727,0.4763186895251477,0.3596254387762725),
    @(23,-1.547866526217289,-1.505981936801861,0.2824917548237356,0.470297959717995,0.3590407684702159),
    @(24,-1.674789248526957,-1.515261708185203,0.2810911706267988,0.4598718992479331,0.3565231970190204),
    @(25,-1.80887009342756,-1.519719238190314,0.2817731283685904,0.4628821637453638,0.3563924736662875),
    @(26,-1.959006920876836,-1.528575197003926,0.283096887754046,0.468178751116758,0.3569955320223329),
    @(27,-2.628002259608796,-1.556738706060644,0.2797088712910335,0.4619532977942937,0.3538870764782091),
    @(28,-2.74519343513492,-1.562936471518154,0.2829719137906281,0.4693568661235901,0.3539594814604756),
    @(29,-2.616699140309548,-1.555467643488537,0.2817799361924212,0.4678291031621125,0.354810156502578),
    @(30,-2.791589609579778,-1.566469941565195,0.2797834832899018,0.4643306684375735,0.354365384517229),
    @(31,-3.652310185618549,-1.60261055199894,0.2759074987536526,0.4603970618121689,0.3493041882624723),
    @(32,-3.756174605209603,-1.606526802227228,0.2739302211430541,0.4620109193277304,0.3478127771358528),
    @(33,-2.882811803245962,-1.568264211766325,0.2801776300458057,0.4677487387920243,0.352676499971505),
    @(34,-3.528348743073225,-1.597079792213081,0.281630330223537,0.4703283830767131,0.3537166474006844),
    @(35,-4.244336103330546,-1.629060565180591,0.2744093277546549,0.4653617167751883,0.3472367264797188),
    @(36,-4.071438086209068,-1.621458705173424,0.2742760764748384,0.4607637103988677,0.3465295138439479),
    @(37,-4.15268585930187,-1.623869896395521,0.2785705132622892,0.4633721880379841,0.3474495646371155),
    @(38,-4.00014195331548,-1.616895041109092,0.2782202908540612,0.4726757433231457,0.3518676810744165),
    @(39,-5.025020935305634,-1.663951270060294,0.2764851337111826,0.4634239827650745,0.3434005295085088),
    @(40,-5.250027027656023,-1.671895129287924,0.2698112294118212,0.4640487639547825,0.3419727282231946),
    @(41,-4.610919432025785,-1.645522452663342,0.2751076341681803,0.4628541490318563,0.3473471770533686)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 14
$lastRow = $startRow + $data.Count - 1

# Copy the formatting of the last existing data row down onto the new rows first.
$srcRange = $ws.Range("A13:F13")
$destRange = $ws.Range("A14:F" + $lastRow)
$srcRange.Copy() | Out-Null
$destRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
